$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> column letter -> new value updates derived from the commit diff
$updates = @(
    @{Row=2; D='28.619.72'; E='  +1.64%  '}
    @{Row=3; D='1.832.11'; E='  +1.56%  '}
    @{Row=4; D='1.001'; E='  -0.03%  '}
    @{Row=5; E='  +0.27%  '}
    @{Row=6; D='1.001'; E='  +0.04%  '}
    @{Row=7; D='0.5346'; E='  +0.91%  '}
    @{Row=8; D='0.3993'; E='  +5.84%  '}
    @{Row=9; D='0.07731'; E='  +3.39%  '}
    @{Row=10; E='  +2.13%  '}
    @{Row=11; D='41.97'; E='  -0.14%  '}
    @{Row=12; D='6.349'; E='  +2.17%  '}
    @{Row=13; D='21.00'; E='  +2.24%  '}
    @{Row=14; D='7.606'; E='  +3.45%  '}
    @{Row=15; D='1.001'; E='  -0.05%  '}
    @{Row=16; D='1.829.31'; E='  +1.92%  '}
    @{Row=17; D='93.57'; E='  +4.45%  '}
    @{Row=18; D='0.00001089'; E='  +2.06%  '}
    @{Row=19; D='0.06571'; E='  +1.02%  '}
    @{Row=20; D='17.78'; E='  +3.13%  '}
    @{Row=21; D='1.001'; E='  +0.09%  '}
    @{Row=22; D='6.089'; E='  +2.85%  '}
    @{Row=23; D='28.632.57'; E='  +1.54%  '}
    @{Row=24; D='11.24'; E='  +0.75%  '}
    @{Row=25; D='2.241'; E='  +7.17%  '}
    @{Row=26; D='20.76'; E='  +1.41%  '}
    @{Row=27; D='2.042.31'; E='  +1.77%  '}
    @{Row=28; D='156.32'; E='  +0.29%  '}
    @{Row=29; D='2.421'; E='  +3.86%  '}
    @{Row=30; D='125.18'; E='  +2.55%  '}
    @{Row=31; D='1.141'; E='  +1.72%  '}
    @{Row=32; D='0.1119'; E='  +2.11%  '}
    @{Row=33; E='  +2.94%  '}
    @{Row=34; D='3.652'}
    @{Row=35; D='0.07251'; E='  +0.30%  '}
    @{Row=36; D='0.2259'; E='  +1.53%  '}
    @{Row=37; E='  +2.44%  '}
    @{Row=38; D='8.899'; E='  +4.86%  '}
    @{Row=39; E='  +2.17%  '}
    @{Row=40; E='  +2.57%  '}
    @{Row=41; D='0.6315'; E='  +2.46%  '}
    @{Row=42; E='  +1.45%  '}
    @{Row=43; D='1.001'; E='  +0.03%  '}
    @{Row=44; E='  -2.69%  '}
    @{Row=45; D='13.58'; E='  +1.29%  '}
    @{Row=46; B='Decentraland'; C='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D='0.5912'; E='  +2.57%  '}
    @{Row=47; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='3.720'; E='  +0.99%  '}
    @{Row=48; D='125.20'; E='  -0.39%  '}
    @{Row=49; D='2.000'; E='  +4.01%  '}
    @{Row=50; D='1.197'; E='  +0.38%  '}
    @{Row=51; D='0.06942'; E='  +1.76%  '}
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) {
        $c = $ws.Cells.Item($u.Row, 2)
        $c.NumberFormat = "@"
        $c.Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $c = $ws.Cells.Item($u.Row, 3)
        $c.NumberFormat = "@"
        $c.Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        $c = $ws.Cells.Item($u.Row, 4)
        $c.NumberFormat = "@"
        $c.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $c = $ws.Cells.Item($u.Row, 5)
        $c.NumberFormat = "@"
        $c.Value = $u.E
    }
}
